# Generate Report for handback
# Fill in the "Correspond Handoff Datetime" (column D) and
# "Correspond Handback DateTime" (column G) values for the handback
# report rows on the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-17 16:16:13"
$wsZhCn.Range("G3").Value = "2016-01-17 16:16:59"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-17 16:16:25"
$wsDeDe.Range("G3").Value = "2016-01-17 16:17:19"
